$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column K
$ws.Range("K1").Value = "tau2"

# tau2 values for rows 2 through 134
$values = @(2.0699999999999998,1.53,0.28000000000000003,3.71,33.08,6.05,0.64,0.68,0.4,77.52,12.37,0.67,2.41,0.2,0.52,15.72,0,0.11,2.5,2.0099999999999998,0.75,3.55,68.11,7.95,1.53,1.32,2.79,64.7,2.65,1.19,7.84,0.001,0.26,0.41,0.08,0.11,0.82,2.0099999999999998,1.1599999999999999,1.01,18.440000000000001,1.05,0.14000000000000001,1.37,3.14,14.22,0.88,0.17,0.63,1.1299999999999999,2.4,2.54,0,16.72,0.41,0.46,3.64,0,21.92,0.16,0,0,0.17,0.18,0.21,0,1.44,3.32,1.66,0.3,0,1.98,6.88,0.7,0.44,0,1.42,1.91,2.46,0.002,0,0.08,0.43,0.55000000000000004,0.48,0,0,0.21,1.92,2.64,0,2,1.19,3.14,1.1299999999999999,32.61,0,0.57999999999999996,1.46,1.19,1.22,20.76,0.7,4.17,0.68,69.25,0.7,0,1.1299999999999999,0.2,0,0.91,0.28999999999999998,0.44,3.51,0.81,7.52,2.19,1.17,38.130000000000003,1.87,1.48,1.17,0.56000000000000005,5.0599999999999996,2.38,3.48,0,3.46,0.81,2.31,1.52,1.64)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Update selection to match the final state (K9 selected, scrolled to top)
$ws.Range("K9").Select()
